$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 27
$ws.Range("D27").Value = "새로운 루다를 지탱하는 모델 서빙 아키텍처 — 2편: ArgoCD와 모델 서빙"
$ws.Range("E27").Value = "https://tech.scatterlab.co.kr/serving-architecture-2/"

# Row 44
$ws.Range("D44").Value = "스타트업 리뷰 (망고부스트)"
$ws.Range("E44").Value = "https://engineering-ladder.tistory.com/118"

# Row 50
$ws.Range("D50").Value = "양자 컴퓨팅"
$ws.Range("E50").Value = "http://incredible.egloos.com/7579491"

# Row 51
$ws.Range("D51").Value = "[Flask] Flask API 서버에 인증 기능 넣기"
$ws.Range("E51").Value = "https://bskyvision.com/entry/Flask-Flask-API-%EC%84%9C%EB%B2%84%EC%97%90-%EC%9D%B8%EC%A6%9D-%EA%B8%B0%EB%8A%A5-%EB%84%A3%EA%B8%B0"

$wb.Save()
